$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values, regenerated to use K instead of Strike# for rows 2-36.
$kValues = @{
    2  = 4
    3  = 1
    4  = 4
    5  = 3
    6  = 5
    7  = 4
    8  = 4
    9  = 7
    10 = 2
    11 = 3
    12 = 1
    13 = 5
    14 = 5
    15 = 1
    16 = 4
    17 = 4
    18 = 2
    19 = 6
    20 = 6
    21 = 3
    22 = 5
    23 = 5
    24 = 1
    25 = 4
    26 = 4
    27 = 4
    28 = 4
    29 = 2
    30 = 4
    31 = 3
    32 = 1
    33 = 5
    34 = 4
    35 = 3
    36 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
